$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test case block mirroring the existing "SampleTest" block (rows 1-2),
# added at rows 5-6 for "ValidateSpeechApi"
$ws.Range("A5").Value = "ValidateSpeechApi"
$ws.Range("B5").Value = "User"
$ws.Range("C5").Value = "Password"
$ws.Range("B6").Value = "sampleusername"
$ws.Range("C6").Value = "samplePwd"

# Match styling of the corresponding Username/Password row (row 2)
$ws.Range("B6:C6").Style = "Hyperlink"

# Update selection to reflect the newly entered range
$ws.Range("B5:C6").Select()
